# Update automatico via Actualizar 03-06-2021 12-08-17
# Shifts the "Ultimo" (last-checked) timestamps in column D down one slot
# and stamps the newest check time into the first block of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newestCheck = 44261.50553879294
$midCheck    = 44261.48421208333
$oldCheck    = 44261.4628730787

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value2 = $newestCheck
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value2 = $midCheck
}

for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value2 = $oldCheck
}
